# Apply the new model ordering in column A (rows 2-26) and set the
# uniform metric values (columns B-I) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "model_11_4_0",
    "model_11_4_22",
    "model_11_4_21",
    "model_11_4_20",
    "model_11_4_19",
    "model_11_4_18",
    "model_11_4_17",
    "model_11_4_16",
    "model_11_4_15",
    "model_11_4_14",
    "model_11_4_13",
    "model_11_4_23",
    "model_11_4_12",
    "model_11_4_10",
    "model_11_4_9",
    "model_11_4_8",
    "model_11_4_7",
    "model_11_4_6",
    "model_11_4_5",
    "model_11_4_4",
    "model_11_4_3",
    "model_11_4_2",
    "model_11_4_1",
    "model_11_4_11",
    "model_11_4_24"
)

$values = @(
    0.3494677884409869,
    0.1470374151965793,
    -0.8846158281072964,
    0.1124530814274172,
    0.7199474573135376,
    1.177122116088867,
    0.2151760458946228,
    0.7244416475296021
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
